$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "0.06106854602694511"
$ws.Cells.Item(2, 2).Value = "0.9811328649520874"
$ws.Cells.Item(2, 3).Value = "0.01186467427760363"
$ws.Cells.Item(2, 4).Value = "0.9981154799461365"

$ws.Cells.Item(3, 1).Value = "0.01020237058401108"
$ws.Cells.Item(3, 2).Value = "0.998121440410614"
$ws.Cells.Item(3, 3).Value = "0.008072340860962868"
$ws.Cells.Item(3, 4).Value = "0.9982501268386841"

$ws.Cells.Item(4, 1).Value = "0.005957551766186953"
$ws.Cells.Item(4, 2).Value = "0.9985502362251282"
$ws.Cells.Item(4, 3).Value = "0.003101204754784703"
$ws.Cells.Item(4, 4).Value = "0.9991923570632935"

$ws.Cells.Item(5, 1).Value = "0.00222117337398231"
$ws.Cells.Item(5, 2).Value = "0.9993669986724854"
$ws.Cells.Item(5, 3).Value = "0.00223432038910687"
$ws.Cells.Item(5, 4).Value = "0.9991923570632935"

$ws.Cells.Item(6, 1).Value = "0.0019832793623209"
$ws.Cells.Item(6, 2).Value = "0.9994282722473145"
$ws.Cells.Item(6, 3).Value = "0.002238910878077149"
$ws.Cells.Item(6, 4).Value = "0.9993269443511963"

$ws.Cells.Item(7, 1).Value = "0.002052647760137916"
$ws.Cells.Item(7, 2).Value = "0.9995507597923279"
$ws.Cells.Item(7, 3).Value = "0.000468796119093895"
$ws.Cells.Item(7, 4).Value = "0.9995961785316467"

$ws.Cells.Item(8, 1).Value = "0.001596540445461869"
$ws.Cells.Item(8, 2).Value = "0.9996733069419861"
$ws.Cells.Item(8, 3).Value = "0.0005505777662619948"
$ws.Cells.Item(8, 4).Value = "0.9997307658195496"

$ws.Cells.Item(9, 1).Value = "0.0009262675885111094"
$ws.Cells.Item(9, 2).Value = "0.9997754096984863"
$ws.Cells.Item(9, 3).Value = "0.0003883853496517986"
$ws.Cells.Item(9, 4).Value = "0.9997307658195496"

$ws.Cells.Item(10, 1).Value = "0.001813531736843288"
$ws.Cells.Item(10, 2).Value = "0.9995712041854858"
$ws.Cells.Item(10, 3).Value = "0.0002218400913989171"
$ws.Cells.Item(10, 4).Value = "1"

$ws.Cells.Item(11, 1).Value = "0.001062522758729756"
$ws.Cells.Item(11, 2).Value = "0.9997345805168152"
$ws.Cells.Item(11, 3).Value = "0.0002457831287756562"
$ws.Cells.Item(11, 4).Value = "0.9998654127120972"

$ws.Cells.Item(12, 1).Value = "0.0007396474247798324"
$ws.Cells.Item(12, 2).Value = "0.9997754096984863"
$ws.Cells.Item(12, 3).Value = "0.00009650475112721324"
$ws.Cells.Item(12, 4).Value = "1"

$ws.Cells.Item(13, 1).Value = "0.0007843737257644534"
$ws.Cells.Item(13, 2).Value = "0.9997549653053284"
$ws.Cells.Item(13, 3).Value = "0.000221568247070536"
$ws.Cells.Item(13, 4).Value = "0.9998654127120972"

$ws.Cells.Item(14, 1).Value = "0.0004416233859956264"
$ws.Cells.Item(14, 2).Value = "0.9998775124549866"
$ws.Cells.Item(14, 3).Value = "0.000126576327602379"
$ws.Cells.Item(14, 4).Value = "1"

$ws.Cells.Item(15, 1).Value = "0.0002139538555638865"
$ws.Cells.Item(15, 2).Value = "0.9999183416366577"
$ws.Cells.Item(15, 3).Value = "0.0001384565985063091"
$ws.Cells.Item(15, 4).Value = "0.9998654127120972"

$ws.Cells.Item(16, 1).Value = "0.001017906935885549"
$ws.Cells.Item(16, 2).Value = "0.9997957944869995"
$ws.Cells.Item(16, 3).Value = "0.00001710674223431852"
$ws.Cells.Item(16, 4).Value = "1"

$ws.Cells.Item(17, 1).Value = "0.000304548884741962"
$ws.Cells.Item(17, 2).Value = "0.9999387264251709"
$ws.Cells.Item(17, 3).Value = "0.00000499692396260798"
$ws.Cells.Item(17, 4).Value = "1"

$ws.Cells.Item(18, 1).Value = "0.0006463424651883543"
$ws.Cells.Item(18, 2).Value = "0.9998366236686707"
$ws.Cells.Item(18, 3).Value = "0.00000937754703045357"
$ws.Cells.Item(18, 4).Value = "1"

$ws.Cells.Item(19, 1).Value = "0.0002065370354102924"
$ws.Cells.Item(19, 2).Value = "0.9999591708183289"
$ws.Cells.Item(19, 3).Value = "0.0005663647898472846"
$ws.Cells.Item(19, 4).Value = "0.9997307658195496"

$ws.Cells.Item(20, 1).Value = "0.0008090437622740865"
$ws.Cells.Item(20, 2).Value = "0.9998570680618286"
$ws.Cells.Item(20, 3).Value = "0.0004376985889393836"
$ws.Cells.Item(20, 4).Value = "0.9997307658195496"

$ws.Cells.Item(21, 1).Value = "0.0004113477189093828"
$ws.Cells.Item(21, 2).Value = "0.9998570680618286"
$ws.Cells.Item(21, 3).Value = "0.0003303957055322826"
$ws.Cells.Item(21, 4).Value = "0.9998654127120972"

$ws.Cells.Item(22, 1).Value = "0.00002181418312829919"
$ws.Cells.Item(22, 2).Value = "1"
$ws.Cells.Item(22, 3).Value = "0.00043743837159127"
$ws.Cells.Item(22, 4).Value = "0.9997307658195496"

$ws.Cells.Item(23, 1).Value = "0.000405898317694664"
$ws.Cells.Item(23, 2).Value = "0.9998570680618286"
$ws.Cells.Item(23, 3).Value = "0.0002248123346362263"
$ws.Cells.Item(23, 4).Value = "0.9998654127120972"

$ws.Cells.Item(24, 1).Value = "0.0004354831180535257"
$ws.Cells.Item(24, 2).Value = "0.9998775124549866"
$ws.Cells.Item(24, 3).Value = "0.001214518211781979"
$ws.Cells.Item(24, 4).Value = "0.9997307658195496"

$ws.Cells.Item(25, 1).Value = "0.0005512360366992652"
$ws.Cells.Item(25, 2).Value = "0.9997549653053284"
$ws.Cells.Item(25, 3).Value = "0.0001320185983786359"
$ws.Cells.Item(25, 4).Value = "0.9998654127120972"

$ws.Cells.Item(26, 1).Value = "0.0001323652977589518"
$ws.Cells.Item(26, 2).Value = "0.999979555606842"
$ws.Cells.Item(26, 3).Value = "0.001012144843116403"
$ws.Cells.Item(26, 4).Value = "0.9997307658195496"

$ws.Cells.Item(27, 1).Value = "0.0002342577063245699"
$ws.Cells.Item(27, 2).Value = "0.999979555606842"
$ws.Cells.Item(27, 3).Value = "0.0008358151535503566"
$ws.Cells.Item(27, 4).Value = "0.9998654127120972"

$ws.Cells.Item(28, 1).Value = "0.0003538480377756059"
$ws.Cells.Item(28, 2).Value = "0.9999183416366577"
$ws.Cells.Item(28, 3).Value = "0.0002917610399890691"
$ws.Cells.Item(28, 4).Value = "0.9998654127120972"

$ws.Cells.Item(29, 1).Value = "0.0004108251596335322"
$ws.Cells.Item(29, 2).Value = "0.9999183416366577"
$ws.Cells.Item(29, 3).Value = "0.001576379174366593"
$ws.Cells.Item(29, 4).Value = "0.9997307658195496"

$ws.Cells.Item(30, 1).Value = "0.0000575275880692061"
$ws.Cells.Item(30, 2).Value = "0.999979555606842"
$ws.Cells.Item(30, 3).Value = "0.003649125341325998"
$ws.Cells.Item(30, 4).Value = "0.9997307658195496"

$ws.Cells.Item(31, 1).Value = "0.000136215821839869"
$ws.Cells.Item(31, 2).Value = "0.9999387264251709"
$ws.Cells.Item(31, 3).Value = "0.00008505047298967838"
$ws.Cells.Item(31, 4).Value = "1"

$ws.Cells.Item(32, 1).Value = "0.0001649706391617656"
$ws.Cells.Item(32, 2).Value = "0.9999387264251709"
$ws.Cells.Item(32, 3).Value = "0.001472195377573371"
$ws.Cells.Item(32, 4).Value = "0.9997307658195496"

$ws.Cells.Item(33, 1).Value = "0.000499738089274615"
$ws.Cells.Item(33, 2).Value = "0.9998978972434998"
$ws.Cells.Item(33, 3).Value = "0.0001446873357053846"
$ws.Cells.Item(33, 4).Value = "1"

$ws.Cells.Item(34, 1).Value = "0.0002880675892811269"
$ws.Cells.Item(34, 2).Value = "0.9998978972434998"
$ws.Cells.Item(34, 3).Value = "0.001405739807523787"
$ws.Cells.Item(34, 4).Value = "0.9997307658195496"

$ws.Cells.Item(35, 1).Value = "0.00005602732562692836"
$ws.Cells.Item(35, 2).Value = "0.999979555606842"
$ws.Cells.Item(35, 3).Value = "0.003235064214095473"
$ws.Cells.Item(35, 4).Value = "0.9997307658195496"

$ws.Cells.Item(36, 1).Value = "0.00008298282773466781"
$ws.Cells.Item(36, 2).Value = "0.9999591708183289"
$ws.Cells.Item(36, 3).Value = "0.003270683577284217"
$ws.Cells.Item(36, 4).Value = "0.9997307658195496"

$ws.Cells.Item(37, 1).Value = "0.0002091178466798738"
$ws.Cells.Item(37, 2).Value = "0.999979555606842"
$ws.Cells.Item(37, 3).Value = "0.00365281687118113"
$ws.Cells.Item(37, 4).Value = "0.9997307658195496"

$ws.Cells.Item(38, 1).Value = "0.0000005787932195744361"
$ws.Cells.Item(38, 2).Value = "1"
$ws.Cells.Item(38, 3).Value = "0.003618230810388923"
$ws.Cells.Item(38, 4).Value = "0.9997307658195496"

$ws.Cells.Item(39, 1).Value = "0.001060109469108284"
$ws.Cells.Item(39, 2).Value = "0.9997957944869995"
$ws.Cells.Item(39, 3).Value = "0.0003030379593838006"
$ws.Cells.Item(39, 4).Value = "0.9998654127120972"

$ws.Cells.Item(40, 1).Value = "0.0000168683964147931"
$ws.Cells.Item(40, 2).Value = "1"
$ws.Cells.Item(40, 3).Value = "0.0009281523525714874"
$ws.Cells.Item(40, 4).Value = "0.9997307658195496"

$ws.Cells.Item(41, 1).Value = "0.0002175613917643204"
$ws.Cells.Item(41, 2).Value = "0.9999387264251709"
$ws.Cells.Item(41, 3).Value = "0.0000003142092452890211"
$ws.Cells.Item(41, 4).Value = "1"

$ws.Cells.Item(42, 1).Value = "0.00003161748463753611"
$ws.Cells.Item(42, 2).Value = "1"
$ws.Cells.Item(42, 3).Value = "0.0000001378857206191242"
$ws.Cells.Item(42, 4).Value = "1"

$ws.Cells.Item(43, 1).Value = "0.0000237984841078287"
$ws.Cells.Item(43, 2).Value = "0.999979555606842"
$ws.Cells.Item(43, 3).Value = "0.00000001059054888230548"
$ws.Cells.Item(43, 4).Value = "1"

$ws.Cells.Item(44, 1).Value = "0.0004486858379095793"
$ws.Cells.Item(44, 2).Value = "0.9998775124549866"
$ws.Cells.Item(44, 3).Value = "0.000001304711076954845"
$ws.Cells.Item(44, 4).Value = "1"

$ws.Cells.Item(45, 1).Value = "0.00006916253187227994"
$ws.Cells.Item(45, 2).Value = "0.9999591708183289"
$ws.Cells.Item(45, 3).Value = "0.001131820026785135"
$ws.Cells.Item(45, 4).Value = "0.9997307658195496"

$ws.Cells.Item(46, 1).Value = "0.0002833192411344498"
$ws.Cells.Item(46, 2).Value = "0.9999387264251709"
$ws.Cells.Item(46, 3).Value = "0.00005508196773007512"
$ws.Cells.Item(46, 4).Value = "1"

$ws.Cells.Item(47, 1).Value = "0.000003831179128610529"
$ws.Cells.Item(47, 2).Value = "1"
$ws.Cells.Item(47, 3).Value = "0.0000004535919515546993"
$ws.Cells.Item(47, 4).Value = "1"

$ws.Cells.Item(48, 1).Value = "0.000001744376504575484"
$ws.Cells.Item(48, 2).Value = "1"
$ws.Cells.Item(48, 3).Value = "0.00000009169063019953683"
$ws.Cells.Item(48, 4).Value = "1"

$ws.Cells.Item(49, 1).Value = "0.000001436220713912917"
$ws.Cells.Item(49, 2).Value = "1"
$ws.Cells.Item(49, 3).Value = "0.00000004345042370346164"
$ws.Cells.Item(49, 4).Value = "1"

$ws.Cells.Item(50, 1).Value = "0.0005121473222970963"
$ws.Cells.Item(50, 2).Value = "0.9998570680618286"
$ws.Cells.Item(50, 3).Value = "0.001869646366685629"
$ws.Cells.Item(50, 4).Value = "0.9997307658195496"

$ws.Cells.Item(51, 1).Value = "0.0001955262705450878"
$ws.Cells.Item(51, 2).Value = "0.9999591708183289"
$ws.Cells.Item(51, 3).Value = "0.00194473983719945"
$ws.Cells.Item(51, 4).Value = "0.9997307658195496"
